# Case and Fatality Demographics Data Updated
# Applies the 2021-03-26 data refresh: new counts for every demographic
# breakdown (age group / gender / race-ethnicity) for both Cases and
# Fatalities, with the "Total"/"Grand Total" cells converted from a
# hard-coded number to a live =SUM() formula.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Cases by Age Group
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value = 269
$ws.Range("B3").Value = 1290
$ws.Range("B4").Value = 3501
$ws.Range("B5").Value = 15292
$ws.Range("B6").Value = 16883
$ws.Range("B7").Value = 14783
$ws.Range("B8").Value = 12408
$ws.Range("B9").Value = 4464
$ws.Range("B10").Value = 3005
$ws.Range("B11").Value = 1788
$ws.Range("B12").Value = 1170
$ws.Range("B13").Value = 1846
$ws.Range("B15").Formula = "=SUM(B2:B14)"
$ws.Activate()
$ws.Range("A20").Select()

# ---------------------------------------------------------------------
# Cases by Gender
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 25732
$ws.Range("B3").Value = 50079
$ws.Range("B4").Value = 901
$ws.Range("B5").Formula = "=SUM(B2:B4)"
$ws.Activate()
$ws.Range("B2:B4").Select()

# ---------------------------------------------------------------------
# Cases by RaceEthnicity
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 942
$ws.Range("B3").Value = 12850
$ws.Range("B4").Value = 27870
$ws.Range("B5").Value = 466
$ws.Range("B6").Value = 26140
$ws.Range("B7").Value = 8444
$ws.Range("B8").Formula = "=SUM(B2:B7)"
$ws.Activate()
$ws.Range("A12").Select()

# ---------------------------------------------------------------------
# Fatalities by Age Group
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B4").Value = 32
$ws.Range("B5").Value = 227
$ws.Range("B6").Value = 760
$ws.Range("B7").Value = 2219
$ws.Range("B8").Value = 5154
$ws.Range("B9").Value = 4375
$ws.Range("B10").Value = 5649
$ws.Range("B11").Value = 6295
$ws.Range("B12").Value = 6293
$ws.Range("B13").Value = 15965
$ws.Range("B15").Formula = "=SUM(B2:B14)"
$ws.Activate()
$ws.Range("C21").Select()

# ---------------------------------------------------------------------
# Fatalities by Gender
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 19665
$ws.Range("B3").Value = 27320
$ws.Range("B5").Formula = "=SUM(B2:B4)"
$ws.Activate()
$ws.Range("C17").Select()

# ---------------------------------------------------------------------
# Fatalities by Race-Ethnicity
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 977
$ws.Range("B3").Value = 4635
$ws.Range("B4").Value = 21773
$ws.Range("B5").Value = 255
$ws.Range("B6").Value = 19321
$ws.Range("B7").Value = 25
$ws.Range("B8").Formula = "=SUM(B2:B7)"
$ws.Activate()
$ws.Range("B12").Select()

# Leave the workbook with "Cases by Age Group" as the active/front sheet,
# matching the saved selection state in the updated file.
$wb.Worksheets.Item("Cases by Age Group").Activate()
